# DW+P5 - Modele plan tests acceptation
# Fill in the Action / Resultat attendu / Resultat observe columns for the
# "produit", "panier" and "confirmation" test rows (rows 3-5), and clear out
# the extra pre-formatted blank rows (6-22) that are no longer needed now
# that the test plan has real content (commit: "Validation of the order
# with empty cart and filled form: fixed").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - "produit" page test
$ws.Range("C3").Value = "Cliquer sur différents produits. Ouvrir les menus déroulants / Taper des données. Cliquer sur le bouton d'ajout au panier"
$ws.Range("D3").Value = "Voire une liste de couleur, pouvoir renseigner une quantité d'article. Possibilité de l'ajouter au panier "
$ws.Range("E3").Value = "OK"

# Row 4 - "panier" page test
$ws.Range("C4").Value = "Modifier les quantités d'un produit, le supprimer. Rentrer les informations adéquates par le formulaires"
$ws.Range("D4").Value = "La quantité totald d'articles ainsi que le prix total doivent être automatiquement calculés et affichés, au gré des modifications et suppressions. Le formulaire doit avertir si les informations entrées sont inadéquats."
$ws.Range("E4").Value = "OK"

# Row 5 - "confirmation" page test
$ws.Range("C5").Value = "Analyser la page pour vérifier si le numéro de commande est bien renseigné"
$ws.Range("D5").Value = "Avoir le numéro de commande bien renseigné"
$ws.Range("E5").Value = "Ok"

# Rows 6-22 were previously pre-formatted (bordered/styled) empty rows left
# over from the template; clear all their content + formatting now that the
# plan only needs rows 1-5, leaving plain rows with just their height.
$ws.Range("A6:E22").Clear()
$ws.Range("A6:A22").EntireRow.RowHeight = 15.75

# Scroll the sheet so row 14 is at the top-left when the file is reopened.
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
